$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 2.316462327490609
$ws.Range("E2").Value = -12.90434384910116

$ws.Range("C3").Value = 9.162319139560072
$ws.Range("E3").Value = 11.00161740362515

$ws.Range("C4").Value = 4.440675708811326
$ws.Range("E4").Value = 0.02592873473603596

$ws.Range("C5").Value = 6.166830030198267
$ws.Range("E5").Value = 5.009734572843771

$ws.Range("C6").Value = 2.204591531842581
$ws.Range("E6").Value = -4.726949348040732

$ws.Range("C7").Value = 1.09689030557385
$ws.Range("E7").Value = -4.442195584720931

$ws.Range("C8").Value = 1.600060471414833
$ws.Range("E8").Value = -2.949819494134909

$ws.Range("C9").Value = 1.651055586686678
$ws.Range("E9").Value = 2.833734372666652

$ws.Range("C10").Value = 1.99440460461342
$ws.Range("E10").Value = 0.9898525518331924

$ws.Range("C11").Value = 2.613530175870626
$ws.Range("E11").Value = 3.313641510414356

$ws.Range("C12").Value = 3.696331036365752
$ws.Range("E12").Value = 6.233648892987009

$ws.Range("C13").Value = -0.973865931199458
$ws.Range("E13").Value = -4.327930935899992

$ws.Range("C14").Value = 1.725130460355095
$ws.Range("E14").Value = -1.194610791900008

$ws.Range("C15").Value = -0.5115004854862049
$ws.Range("E15").Value = 2.23509962177757

$ws.Range("C16").Value = 0.9727820482463123
$ws.Range("E16").Value = 2.866869504079239

$ws.Range("C17").Value = 2.122313752051319
$ws.Range("E17").Value = 0.8382457967197388

$ws.Range("C18").Value = -0.1133034947815914
$ws.Range("E18").Value = 0.8323378752418176

$ws.Range("C19").Value = 3.161594928268019
$ws.Range("E19").Value = 2.57979941834241
